# Actualizacion inter_cabina_casos de prueba 13/07/2015
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = "CU1, CU2, CU3`n"
$ws.Range("C12").Value = "CU4"
$ws.Range("C13").Value = "CU5, CU6`n"
$ws.Range("C14").Value = "CU5, CU7"
$ws.Range("C15").Value = "CU5, CU8"
$ws.Range("C16").Value = "CU8, CU9, CU12, CU13, CU15, CU14"
$ws.Range("C17").Value = "CU9, CU10"
$ws.Range("C18").Value = "CU11"
$ws.Range("C19").Value = "CU14"
$ws.Range("C20").Value = "N/A"
$ws.Range("C21").Value = "N/A"
$ws.Range("C22").Value = "N/A"
$ws.Range("C23").Value = "N/A"

$ws.Range("D23").Select()
